$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unmerge the old multi-level header merged ranges ---
$ws.Range("H1:L1").UnMerge()
$ws.Range("M1:P1").UnMerge()
$ws.Range("Q1:S1").UnMerge()

# --- Row 1: flatten the two-row header into one descriptive header row ---
$headerVals = @("Player ID","Player","#","Nation","Pos","Age","90s","Tkl","TklW","Def 3rd","Mid 3rd","Att 3rd","Cha","Att","Tkl%","Lost","Blocks","Sh","Pass","Int","Tkl+Int","Clr","Err")
for ($i = 0; $i -lt $headerVals.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headerVals[$i]
}

# --- Row 2 keeps the old sub-header labels, but the row itself is now hidden ---
$ws.Rows.Item(2).Hidden = $true

# --- Row 3 is a blank spacer row, now hidden ---
$ws.Rows.Item(3).Hidden = $true

# --- Fill in previously-blank "Blocks" values with 0 ---
$blocksRows = @(4, 6, 10, 11, 17, 19)
foreach ($r in $blocksRows) {
  $ws.Cells.Item($r, 15).Value = 0
}

# --- Row 20 (the summary row) is now hidden ---
$ws.Rows.Item(20).Hidden = $true

# --- Update the active selection to mirror the saved workbook state ---
$ws.Range("O21").Select()
